$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.960.93'
$ws.Range("E2").Value = '  -1.27%  '
$ws.Range("D3").Value = '3.146.24'
$ws.Range("E3").Value = '  -0.61%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '603.32'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.18%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.32'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.84%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").Value = '3.139.71'
$ws.Range("E8").Value = '  -0.77%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.528'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.52%  '
$ws.Range("E10").Value = '  -1.79%  '
$ws.Range("E11").Value = '  -1.96%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.468'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.44%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000254'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.54%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.09'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.23%  '
$ws.Range("D15").Value = '3.666.75'
$ws.Range("E15").Value = '  -0.51%  '
$ws.Range("E16").Value = '  +2.06%  '
$ws.Range("D17").Value = '64.008.15'
$ws.Range("E17").Value = '  -1.16%  '
$ws.Range("D18").Value = '3.141.47'
$ws.Range("E18").Value = '  -0.74%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.88'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.78%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '488.51'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.72%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.73'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.14%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.714'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.95%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.75'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.89%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '88.29'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +4.21%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '13.31'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.53%  '
$ws.Range("E26").Value = '  +0.06%  '
$ws.Range("E27").Value = '  -2.05%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.22'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.70%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.02'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.61%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.07'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.29%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '27.71'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.71%  '
$ws.Range("E32").Value = '  -6.14%  '
$ws.Range("E33").Value = '  +0.00%  '
$ws.Range("E34").Value = '  -1.77%  '
$ws.Range("E35").Value = '  -3.21%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.07'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.55%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '52.66'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.01%  '
$ws.Range("E38").Value = '  -5.47%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.96'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -7.08%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0398'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.79%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '434.16'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -6.86%  '
$ws.Range("E42").Value = '  -0.47%  '
$ws.Range("E43").Value = '  -0.22%  '
$ws.Range("D44").Value = '2.933.79'
$ws.Range("E44").Value = '  +3.00%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.261'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.94%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.20'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -5.89%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.41'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.72%  '
$ws.Range("E48").Value = '  -0.09%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '25.92'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.26%  '
$ws.Range("E50").Value = '  +0.43%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '120.52'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.19%  '
